$d = $word.ActiveDocument

# --- Paragraph 1 (title block): replace both runs' text ---
$d.Content.Find.Execute(
    "🚀המאמר היומי של מייק 19.09.24: ⚡️🚀 ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "🚀המאמר היומי של מייק 17.09.24: ⚡️🚀", 2)

$d.Content.Find.Execute(
    "Training Chain-of-Thought via Latent-Variable Inference",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "STaR: Self-Taught Reasoner Bootstrapping Reasoning With Reasoning", 2)

# --- Paragraph 2 ---
$d.Content.Find.Execute(
    "ממשיכים בקו הסקירות שהובילו (לפחות לעניות דעתי) למודל החדש (יחסית, יצא כבר לפני שבוע) של openai. במאמר הקודם שסקרתי STaR דיברנו על איך ניתן לשפר יכולת ריזונינג של מודל שפה כאשר יש בידינו דאטהסט גדול יחסית של שאלות ותשובות D ודאטהסט קטן של שאלות ותשובות עם הריזונינג. בגדול הרעיון שם היא לרתום מודל שפה לייצר ריזונינג לשאלות, להוסיף שאלות שהריזונינג שלהם הוביל לתשובה נכונה לדאטהסט הקטן ולהמשיך לאמן עד ההתכנסות.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "אני ממשיך לחפור במאמרי שאולי עיצבו את הנתיב הובילו ל-o1 של openai. הפעם נברתי כה עמוק שהגעתי למאמר שיצא לפני שנתיים וחצי (בדיפ היום זה כמו 100 שנה במתמטיקה). שימו לב שהמאמר יצא עוד לפני chatgpt. המאמר הזה מציע שיטה לשיפור יכולת reasoning של מודל שפה כאשר בידנו יש דאטהסט גדול של שאלות ותשובות D ודאטהסט קטן D_R הרבה יותר (המאמר מדבר על 10 דוגמאות בלבד) המכיל בנוסף גם את שרשרת ה-reasoning.", 2)

# --- Paragraph 3 (target text contains two literal straight apostrophes;
#     Find/Replace auto-curls "'" to a right single quote, so drop in a
#     placeholder token here and patch the two apostrophes back in below) ---
$d.Content.Find.Execute(
    "המאמר הנוכחי שיצא בערך שנה וחצי אחריו משכלל את הגישה הזו ומציע שיטה ש״ממנפת״ גם את השאלות שעבוד המודל יצר ריזונינג שלא הוביל לתשובה הנכונה. המאמר מכיל מתמטיקה די כבדה אז אנסה להעביר לכם את הרעיון הכללי יחסית בפשטות. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "כאשר אני מדבר על שיפור איכות ה-reasoning אני בעצם מתכוון לפיינטיון של המודל במטרה לקבל מודל חזק יותר ב-reasoning. המחברים מציעים אלגוריתם המורכב משני שלבים עיקריים. בשלב הראשון מזינים את הבאץZZAPOSZZ של שאלות למודל שפה כאשר בנוסף לשאלות הפרומפט מכיל את דוגמאות לשרשראות ה-reasoning m מ- D_R. המודל מתבקש לבנות שרשרת reasoning לכל השאלות מבאץZZAPOSZZ (לא מ-D_R) ולהגיע לתשובה הסופית.", 2)

for ($i = 0; $i -lt 2; $i++) {
    $r = $d.Content
    $r.Find.Text = "ZZAPOSZZ"
    $r.Find.Forward = $true
    $r.Find.Wrap = 0
    $found = $r.Find.Execute()
    if ($found) {
        $r.Text = "'"
    }
}

# --- Paragraph 4 ---
$d.Content.Find.Execute(
    "הרי המטרה שלנו היא לעשות פיינטיון למודל שפה כך שיכולת הריזונינג שלו תשתפר. מתמטית ניתן לתרגם את הבעיה לבעיה וריאציונית באופן הבא. אנו מעוניינים לאמן מודל שיוצר ריזונינג עבור שאלה x. מה שיש לנו זה דאטהסט של שאולות x ו-תשובות y. אז אנחנו רוצים לאמן את המודל להפיק ריזונינג z (ניתן להתייחס אלי כמו אל משתנה לטנטי) מהתפלגות בהינתן השאלה x מ-D תוך כדי ניצול של התשובה y. כלומר אנו רוצים למקסם את הנראות (likelihood) של ההתפלגות המותנית של הריזונינג z בהינתן (עבור) שאלה x ותשובה y. במילים פשוטות אנו מאפטמים את פרמטרי המודל כך שהנראות הזו תהיה מקסימלית על D.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "את שרשראות ה-reasoning לשאלות שהצליחו להגיע לתשובה נכונה מוסיפים לסט שנקרא לו D_N. לשאלות שהמודל לא הצליח להגיע לתשובה סופית נכונה אנחנו מוסיפים רמז (במאמר זה נקרא rationalization) שעוזר למודל לבנות את שרשרת ה-reasoning. השאלות שהצליחו להגיע לתשובה הנכונה אחרי הרמז גם נוספים ל D_N. לאחר מכן מבצעים איטרציה אחת של שיטת מורד הגרדיאנט נבחרת על D_N ומעדכנים את משקלי המודל. חוזרים על השלבים האלו עד שהלוס מתייצב.", 2)

# --- Paragraph 5 (keeps a trailing space / xml:space=preserve) ---
$d.Content.Find.Execute(
    "אולם אנו לא יכולים לעשות זאת בצורה ישירה כלומר לא ניתן לדגום את הריזונינג בהינתן שאלה x ותשובה y. הסיבה לכך היא שאנו לא רוצים לאמן מודל שמייצר ריזונינג לשאלה יחד עם התשובה (כי אנו רוצים מודל שיפתור לנו שאלות בלי לדעת את התשובה). אז המאמר הקודם בחר לנצל את תשובה y על ידי פלטור החוצה של z שהובילו לתשובות לא נכונות. לעומת זאת המאמר הזה מציע שיטה שבה אנו ממנפים גם את ה- z-ים הלא נכונים לשיפור המודל. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "זהו זה, שיטה אינטואיטיבית ופשוטה שקיבלה כמה מאמרי השמך די כבדים שבתקווה אסקור אותם גם כן ", 2)

# --- Paragraph 6 (becomes the final URL line) ---
$d.Content.Find.Execute(
    "כאמור המאמר מנצל כמה שיטות מתמטיות די כבדות לכך ואחת מהם הוא שכלול של Markov Chain Monte Carlo כאשר ה-proposal distribution  (שממנו דוגמים במטרה שזו תתכנס עם הזמן להתפלגות היעד כלומר זו של ריזונינג z בהינתן שאלה x ותשובה y) משתנה עם האיטרציה להאצת התכנסות (Markovian score climbing שהוא שכלול של Robbins-Monro לחישוב ״גודל העדכון״).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "https://arxiv.org/pdf/2203.14465", 2)

# --- Remove trailing paragraphs 7-10 (now at indices 7-10), whose content
#     ("מה הקשר ל-MCMC...", "בנוסף המאמר משכלל...", "וכל זה כדי לשפר...",
#      old URL "https://arxiv.org/pdf/2312.02179") is dropped entirely. ---
$start = $d.Paragraphs(7).Range.Start
$end = $d.Paragraphs($d.Paragraphs.Count).Range.End
$r = $d.Range($start, $end)
$r.Delete()
